# S21/G01: Custom bracket orders (primary + follow-up GTT legs)
# Appends 4 new task rows (174-177) to the sprint tasks sheet, covering
# S21/G01 (custom bracket orders) and the start of S21/G02 (bracket-order
# backtesting), and extends the sheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A174").Value = "S21"
$ws.Range("B174").Value = "G01"
$ws.Range("C174").Value = "Custom bracket orders (primary + follow-up GTT legs)"
$ws.Range("D174").Value = "S21_G01_TB001"
$ws.Range("E174").Value = "Design and implement backend helper/flow to create paired manual orders (primary + LIMIT GTT leg) given side, qty, effective price, and MTP."
$ws.Range("F174").Value = "Keeps invariant that each leg is a normal WAITING manual order; does not introduce broker-specific bracket types."
$ws.Range("G174").Value = "planned"
$ws.Range("H174").Value = "Provides a single place to compute P_target/P_reentry and create the extra GTT order consistently for both BUY and SELL."

$ws.Range("A175").Value = "S21"
$ws.Range("B175").Value = "G01"
$ws.Range("C175").Value = "Custom bracket orders (primary + follow-up GTT legs)"
$ws.Range("D175").Value = "S21_G01_TF001"
$ws.Range("E175").Value = "Extend Holdings Buy/Sell dialog with a Bracket section that lets the user enable a follow-up GTT leg, pre-fills MTP% from current appreciation, and previews the derived target price."
$ws.Range("F175").Value = "Bracket invocation remains manual; dialog simply creates two manual orders via the existing orders API."
$ws.Range("G175").Value = "planned"
$ws.Range("H175").Value = "Makes it easy to add profit-target or re-entry GTT orders alongside normal trades without changing queue semantics."

$ws.Range("A176").Value = "S21"
$ws.Range("B176").Value = "G01"
$ws.Range("C176").Value = "Custom bracket orders (primary + follow-up GTT legs)"
$ws.Range("D176").Value = "S21_G01_TF002"
$ws.Range("E176").Value = "Highlight bracket-related information in the Queue and Orders grids (order_type, trigger_price, GTT flag) and optionally tag bracket legs for easier identification."
$ws.Range("F176").Value = "Builds on the existing DataGrid-based queue and orders views added earlier."
$ws.Range("G176").Value = "planned"
$ws.Range("H176").Value = "Improves transparency so you can quickly see which waiting orders are bracket legs and how they are configured."

$ws.Range("A177").Value = "S21"
$ws.Range("B177").Value = "G02"
$ws.Range("C177").Value = "Bracket-order backtesting using Kite OHLCV"
$ws.Range("D177").Value = "S21_G02_TB001"
$ws.Range("E177").Value = "Add a backend console script that pulls OHLCV via the existing market-data layer and simulates the custom bracket logic over a given symbol, timeframe, and lookback."
$ws.Range("F177").Value = "Focus initial experiments on BSE and NETWEB over the last month but keep the script parameterised for any symbol."
$ws.Range("G177").Value = "planned"
$ws.Range("H177").Value = "Enables quantitative evaluation of how often the MTP-based bracket legs would have been filled and the resulting P&L and drawdowns."
